$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates, applied identically to the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets which share the same data.
$updates = @{
    2  = 1092
    3  = 799
    4  = 269
    5  = 44
    6  = 1105
    8  = 2033
    9  = 7557
    10 = 898
    11 = 412
    12 = 342
    13 = 127
    16 = 7100
    18 = 1331
    22 = 137
    23 = 297
    24 = 135
    28 = 19
    29 = 408
    34 = 74
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
